$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.018808245658875
$ws.Range("B1").Value = 1.919378042221069
$ws.Range("C1").Value = 2.821298122406006
$ws.Range("D1").Value = 3.443411350250244
$ws.Range("E1").Value = 2.049438714981079
